$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "수학&통계학 for MSDS 난이도 (2)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/boot-camp-difficulty-2/#utm_source=rss&utm_medium=rss&utm_campaign=boot-camp-difficulty-2"

$ws.Range("D36").Value = "Fine-Grained Named Entity Recognition"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/322"

$ws.Range("D37").Value = "[Paper Review] SMART: Robust and Efficient Fine-Tuning for Pre-trained Natural Language Models through Principled Regularized O"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1462&mod=document&pageid=1"

$ws.Range("D46").Value = "[Bioinformatics] 2021년 05월, 의료 인공지능 전문가 양성과정 교육생 모집"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/394"
